$wb = $excel.ActiveWorkbook

# This edit adds a brand-new localization record for the file
# "14205dfa-0802-4929-813f-fa71b3bf7992.md" as the new first data row
# (row 2) on every sheet, pushing the pre-existing "81a5b25e-...md"
# record down to row 3 (with all of its original values, just shifted).

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/764993c81838d7910134f70557271476acfc2b19/e2e/"

# ===================== Sheet: Overview =====================
$wsOv = $wb.Worksheets.Item("Overview")

# Push the existing "81a5b25e" row down from row 2 to row 3.
$wsOv.Range("A3").Value = $wsOv.Range("A2").Value2
$wsOv.Range("B3").Value = $wsOv.Range("B2").Value2
$wsOv.Range("C3").Value = $wsOv.Range("C2").Value2
$wsOv.Range("D3").Value = $wsOv.Range("D2").Value2
$wsOv.Range("E3").Value = $wsOv.Range("E2").Value2
$wsOv.Range("F3").Value = $wsOv.Range("F2").Value2
$wsOv.Range("G3").Value = $wsOv.Range("G2").Value2

# Write the new "14205dfa" record into row 2.
$wsOv.Range("A2").Value = "14205dfa-0802-4929-813f-fa71b3bf7992.md"
$wsOv.Range("B2").Value = "e2e\14205dfa-0802-4929-813f-fa71b3bf7992.md"
$wsOv.Range("C2").Value = ".md"
$wsOv.Range("D2").Value = ""
$wsOv.Range("E2").Value = "Ready for handoff"
$wsOv.Range("F2").Value = "Ready for handoff"
$wsOv.Range("G2").Value = "2016-08-30 18:50:13"

# Rebuild the hyperlinks for column B: new row 2 points at the new file,
# row 3 keeps pointing at the original 81a5b25e file.
$wsOv.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), ($githubBase + "14205dfa-0802-4929-813f-fa71b3bf7992.md"), "", "", "e2e\14205dfa-0802-4929-813f-fa71b3bf7992.md")
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), ($githubBase + "81a5b25e-1357-4f92-ab21-31b8813c4143.md"), "", "", "e2e\81a5b25e-1357-4f92-ab21-31b8813c4143.md")

# Grow the "Overview" table to include the new row.
$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ===================== Sheet: zh-cn =====================
$wsZh = $wb.Worksheets.Item("zh-cn")

# Push the existing "81a5b25e" row down from row 2 to row 3.
$wsZh.Range("A3").Value = $wsZh.Range("A2").Value2
$wsZh.Range("B3").Value = $wsZh.Range("B2").Value2
$wsZh.Range("C3").Value = $wsZh.Range("C2").Value2
$wsZh.Range("D3").Value = $wsZh.Range("D2").Value2
$wsZh.Range("E3").Value = $wsZh.Range("E2").Value2
$wsZh.Range("F3").Value = "'" + $wsZh.Range("F2").Value2
$wsZh.Range("G3").Value = $wsZh.Range("G2").Value2
$wsZh.Range("H3").Value = $wsZh.Range("H2").Value2
$wsZh.Range("I3").Value = $wsZh.Range("I2").Value2
$wsZh.Range("J3").Value = $wsZh.Range("J2").Value2
$wsZh.Range("K3").Value = $wsZh.Range("K2").Value2
$wsZh.Range("L3").Value = $wsZh.Range("L2").Value2
$wsZh.Range("M3").Value = "'" + $wsZh.Range("M2").Value2
$wsZh.Range("N3").Value = $wsZh.Range("N2").Value2
$wsZh.Range("O3").Value = "'" + $wsZh.Range("O2").Value2
$wsZh.Range("P3").Value = $wsZh.Range("P2").Value2

# Write the new "14205dfa" record into row 2.
$wsZh.Range("A2").Value = "14205dfa-0802-4929-813f-fa71b3bf7992.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "'False"
$wsZh.Range("G2").Value = "14205dfa-0802-4929-813f-fa71b3bf7992.68c592ff0fef0d27b6f488348c60488ce11d253a.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-30 18:49:57"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"
$wsZh.Range("L2").Value = ""
$wsZh.Range("M2").Value = "'True"
$wsZh.Range("N2").Value = ""
$wsZh.Range("O2").Value = "'False"
$wsZh.Range("P2").Value = ""

# Rebuild the hyperlinks for column A.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($githubBase + "14205dfa-0802-4929-813f-fa71b3bf7992.md"), "", "", "14205dfa-0802-4929-813f-fa71b3bf7992.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($githubBase + "81a5b25e-1357-4f92-ab21-31b8813c4143.md"), "", "", "81a5b25e-1357-4f92-ab21-31b8813c4143.md")

# Grow the "zh-cn" table to include the new row.
$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ===================== Sheet: de-de =====================
$wsDe = $wb.Worksheets.Item("de-de")

# Push the existing "81a5b25e" row down from row 2 to row 3.
$wsDe.Range("A3").Value = $wsDe.Range("A2").Value2
$wsDe.Range("B3").Value = $wsDe.Range("B2").Value2
$wsDe.Range("C3").Value = $wsDe.Range("C2").Value2
$wsDe.Range("D3").Value = $wsDe.Range("D2").Value2
$wsDe.Range("E3").Value = $wsDe.Range("E2").Value2
$wsDe.Range("F3").Value = "'" + $wsDe.Range("F2").Value2
$wsDe.Range("G3").Value = $wsDe.Range("G2").Value2
$wsDe.Range("H3").Value = $wsDe.Range("H2").Value2
$wsDe.Range("I3").Value = $wsDe.Range("I2").Value2
$wsDe.Range("J3").Value = $wsDe.Range("J2").Value2
$wsDe.Range("K3").Value = $wsDe.Range("K2").Value2
$wsDe.Range("L3").Value = $wsDe.Range("L2").Value2
$wsDe.Range("M3").Value = "'" + $wsDe.Range("M2").Value2
$wsDe.Range("N3").Value = $wsDe.Range("N2").Value2
$wsDe.Range("O3").Value = "'" + $wsDe.Range("O2").Value2
$wsDe.Range("P3").Value = $wsDe.Range("P2").Value2

# Write the new "14205dfa" record into row 2.
$wsDe.Range("A2").Value = "14205dfa-0802-4929-813f-fa71b3bf7992.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "'False"
$wsDe.Range("G2").Value = "14205dfa-0802-4929-813f-fa71b3bf7992.68c592ff0fef0d27b6f488348c60488ce11d253a.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-30 18:50:13"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDe.Range("L2").Value = ""
$wsDe.Range("M2").Value = "'True"
$wsDe.Range("N2").Value = ""
$wsDe.Range("O2").Value = "'False"
$wsDe.Range("P2").Value = ""

# Rebuild the hyperlinks for column A.
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($githubBase + "14205dfa-0802-4929-813f-fa71b3bf7992.md"), "", "", "14205dfa-0802-4929-813f-fa71b3bf7992.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($githubBase + "81a5b25e-1357-4f92-ab21-31b8813c4143.md"), "", "", "81a5b25e-1357-4f92-ab21-31b8813c4143.md")

# Grow the "de-de" table to include the new row.
$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

Write-Host "Generate Report for Handoff: added 14205dfa-0802-4929-813f-fa71b3bf7992.md record to Overview, zh-cn and de-de sheets."
